# rda.xlsx fixture update:
#  - extend the "79.65 Lebenslanges Lernen" classification cell with a second line
#  - add a new "Fingerprint" worksheet (cloned layout/structure of the
#    "Bibliographische Zitate" sheet) describing Pica field 026 $e $5
#  - move the active sheet/tab selection: Basisklassifikation becomes active,
#    with C14 selected; Bibliographische Zitate is no longer the selected tab

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Beschreibung
$ws2 = $wb.Worksheets.Item(2)   # Basisklassifikation
$ws3 = $wb.Worksheets.Item(3)   # Bibliographische Zitate

# --- Basisklassifikation: append a second classification line -------------
$ws2.Range("C14").Value = "79.65 Lebenslanges Lernen`n81.92 Berufliche Weiterbildung"

# --- restore/track the per-sheet selections that existed before our edits -
$ws1.Select()
$ws1.Range("E16").Select()

$ws3.Select()
$ws3.Range("C12").Select()

# --- new "Fingerprint" worksheet, cloned from "Bibliographische Zitate" ---
$ws3.Copy($null, $ws3)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "Fingerprint"

$ws4.Range("B1").Value = "Fingerprint"
$ws4.Range("B2").Value = "Fingerprint bei alten Drucken mit Herkunftsangabe aus Pica-Feld 2275 mit Unterfeld `$a"
$ws4.Range("B3").Value = "Fingerprint Identifier (Unparsed fingerprint, Institution to which field applies)"
$ws4.Range("B5").Value = "Ticket #149"

$ws4.Range("A8").Value = "026 `$e `$5"
$ws4.Range("B8").Value = "Fingerprint"

$ws4.Range("A12").Value = "026 `$e `$5"
$ws4.Range("B12").Value = 151797196
$ws4.Range("C12").Value = "S.ME u-r- m-r- siRe 3 1700R (UFB Erfurt/Gotha; NLB Hannover; Bibliothek des Herzog Anton Ulrich-Museums Braunschweig; SBB-PK Berlin)"

$ws4.Range("A13").Value = "026 `$e `$5"
$ws4.Range("B13").Value = 770927416
$ws4.Range("C13").Value = "t,n, o-s- e-n- Dese C 1539A (SBB; ThULB Jena; Wartburg-Stiftung Eisenach)"

# row heights follow the content now shown (auto-fit approximations)
$ws4.Rows.Item(2).RowHeight = 39.55
$ws4.Rows.Item(12).RowHeight = 39.55
$ws4.Rows.Item(13).RowHeight = 26.85

# hyperlink on B5 should point at ticket #149, not the cloned #150
$ws4.Hyperlinks.Delete()
$ws4.Hyperlinks.Add($ws4.Range("B5"), "http://redmine.thulb.uni-jena.de/issues/149", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "Ticket #149")
# adding the hyperlink auto-applies Excel's default hyperlink character style;
# restore the original (non-hyperlink) cell formatting, same as B3's
$ws4.Range("B3").Copy()
$ws4.Range("B5").PasteSpecial(-4122) # xlPasteFormats

$ws4.Range("C12").Select()

# --- Basisklassifikation becomes the active/selected tab, cell C14 --------
$ws2.Select()
$ws2.Range("C14").Select()
